$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 356.9
$ws.Cells.Item(28, 9).Value = 411.57144
$ws.Cells.Item(28, 11).Value = 411.57144
$ws.Cells.Item(28, 13).Value = 73.42856
$ws.Cells.Item(33, 8).Value = 309.23914
$ws.Cells.Item(33, 9).Value = 279.35715
$ws.Cells.Item(33, 10).Value = 623
$ws.Cells.Item(33, 11).Value = 279.35715
$ws.Cells.Item(33, 12).Value = 623
$ws.Cells.Item(33, 13).Value = -50.35714999999999
$ws.Cells.Item(33, 14).Value = -1081
$ws.Cells.Item(76, 8).Value = 6175750.5
$ws.Cells.Item(76, 9).Value = 3200
$ws.Cells.Item(76, 10).Value = 18520852
$ws.Cells.Item(76, 11).Value = 3200
$ws.Cells.Item(76, 12).Value = 18520852
$ws.Cells.Item(76, 13).Value = -2885
$ws.Cells.Item(76, 14).Value = -18521482
$ws.Cells.Item(79, 8).Value = 6175750.5
$ws.Cells.Item(79, 9).Value = 3200
$ws.Cells.Item(79, 10).Value = 18520852
$ws.Cells.Item(79, 11).Value = 3200
$ws.Cells.Item(79, 12).Value = 18520852
$ws.Cells.Item(79, 13).Value = -2108
$ws.Cells.Item(79, 14).Value = -18523036
$ws.Cells.Item(125, 8).Value = 1157.7142
$ws.Cells.Item(125, 9).Value = 532.5
$ws.Cells.Item(125, 11).Value = 4792.5
$ws.Cells.Item(125, 13).Value = -2332.5
$ws.Cells.Item(132, 8).Value = 60006.11
$ws.Cells.Item(132, 9).Value = 71647.60000000001
$ws.Cells.Item(132, 10).Value = 1798.6666
$ws.Cells.Item(132, 11).Value = 214942.8
$ws.Cells.Item(132, 12).Value = 5395.9998
$ws.Cells.Item(132, 13).Value = -212412.8
$ws.Cells.Item(132, 14).Value = -10455.9998
$ws.Cells.Item(137, 8).Value = 19779.785
$ws.Cells.Item(137, 9).Value = 2025.8182
$ws.Cells.Item(137, 10).Value = 45252.87
$ws.Cells.Item(137, 11).Value = 6077.4546
$ws.Cells.Item(137, 12).Value = 135758.61
$ws.Cells.Item(137, 13).Value = -3527.4546
$ws.Cells.Item(137, 14).Value = -140858.61
$ws.Cells.Item(138, 8).Value = 3273.6667
$ws.Cells.Item(138, 9).Value = 2680.5
$ws.Cells.Item(138, 11).Value = 8041.5
$ws.Cells.Item(138, 13).Value = -2901.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14733.0625
$ws.Cells.Item(32, 9).Value = 15766.172
$ws.Cells.Item(32, 10).Value = 7501.3
$ws.Cells.Item(32, 11).Value = 15766.172
$ws.Cells.Item(32, 12).Value = 7501.3
$ws.Cells.Item(32, 13).Value = -15479.172
$ws.Cells.Item(32, 14).Value = -8075.3
$ws.Cells.Item(45, 8).Value = 3332.842
$ws.Cells.Item(45, 9).Value = 2388.75
$ws.Cells.Item(45, 10).Value = 4019.4546
$ws.Cells.Item(45, 11).Value = 2388.75
$ws.Cells.Item(45, 12).Value = 4019.4546
$ws.Cells.Item(45, 13).Value = -2011.75
$ws.Cells.Item(45, 14).Value = -4773.4546
$ws.Cells.Item(70, 8).Value = 40000
$ws.Cells.Item(70, 10).Value = 40000
$ws.Cells.Item(70, 12).Value = 40000
$ws.Cells.Item(70, 14).Value = -40540
$ws.Cells.Item(73, 8).Value = 40000
$ws.Cells.Item(73, 10).Value = 40000
$ws.Cells.Item(73, 12).Value = 40000
$ws.Cells.Item(73, 14).Value = -41872
$ws.Cells.Item(74, 8).Value = 1450.5625
$ws.Cells.Item(74, 9).Value = 1016.125
$ws.Cells.Item(74, 11).Value = 1016.125
$ws.Cells.Item(74, 13).Value = -142.125
$ws.Cells.Item(77, 8).Value = 1450.5625
$ws.Cells.Item(77, 9).Value = 1016.125
$ws.Cells.Item(77, 11).Value = 5080.625
$ws.Cells.Item(77, 13).Value = -712.625
$ws.Cells.Item(102, 8).Value = 4876.25
$ws.Cells.Item(102, 9).Value = 3168.3333
$ws.Cells.Item(102, 11).Value = 3168.3333
$ws.Cells.Item(102, 13).Value = -1546.3333
$ws.Cells.Item(129, 8).Value = 41417.4
$ws.Cells.Item(129, 9).Value = 7090
$ws.Cells.Item(129, 11).Value = 7090
$ws.Cells.Item(129, 13).Value = -2090
$ws.Cells.Item(132, 8).Value = 26415.857
$ws.Cells.Item(132, 9).Value = 2060.0833
$ws.Cells.Item(132, 11).Value = 6180.249899999999
$ws.Cells.Item(132, 13).Value = -3650.249899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 866.6667
$ws.Cells.Item(20, 9).Value = 866.6667
$ws.Cells.Item(20, 11).Value = 866.6667
$ws.Cells.Item(20, 13).Value = -619.6667
$ws.Cells.Item(86, 8).Value = 1627.091
$ws.Cells.Item(86, 9).Value = 1423.0769
$ws.Cells.Item(86, 10).Value = 1921.7778
$ws.Cells.Item(86, 11).Value = 1423.0769
$ws.Cells.Item(86, 12).Value = 1921.7778
$ws.Cells.Item(86, 13).Value = -300.0769
$ws.Cells.Item(86, 14).Value = -4167.7778
$ws.Cells.Item(89, 8).Value = 1627.091
$ws.Cells.Item(89, 9).Value = 1423.0769
$ws.Cells.Item(89, 10).Value = 1921.7778
$ws.Cells.Item(89, 11).Value = 7115.3845
$ws.Cells.Item(89, 12).Value = 9608.889000000001
$ws.Cells.Item(89, 13).Value = -1499.3845
$ws.Cells.Item(89, 14).Value = -20840.889
$ws.Cells.Item(94, 8).Value = 4499.875
$ws.Cells.Item(94, 9).Value = 1999.75
$ws.Cells.Item(94, 10).Value = 7000
$ws.Cells.Item(94, 11).Value = 1999.75
$ws.Cells.Item(94, 12).Value = 7000
$ws.Cells.Item(94, 13).Value = -1548.75
$ws.Cells.Item(94, 14).Value = -7902
$ws.Cells.Item(99, 8).Value = 2296.5
$ws.Cells.Item(99, 9).Value = 2200
$ws.Cells.Item(99, 11).Value = 2200
$ws.Cells.Item(99, 13).Value = -702
$ws.Cells.Item(105, 9).Value = 1637.1428
$ws.Cells.Item(105, 10).Value = 3847769.2
$ws.Cells.Item(105, 11).Value = 1637.1428
$ws.Cells.Item(105, 12).Value = 3847769.2
$ws.Cells.Item(105, 13).Value = 109.8571999999999
$ws.Cells.Item(105, 14).Value = -3851263.2
$ws.Cells.Item(107, 8).Value = 1186.3158
$ws.Cells.Item(107, 9).Value = 1214.1177
$ws.Cells.Item(107, 10).Value = 950
$ws.Cells.Item(107, 11).Value = 1214.1177
$ws.Cells.Item(107, 12).Value = 950
$ws.Cells.Item(107, 13).Value = 705.8823
$ws.Cells.Item(107, 14).Value = -4790

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 11955.116
$ws.Cells.Item(31, 9).Value = 29729.5
$ws.Cells.Item(31, 10).Value = 3374.3794
$ws.Cells.Item(31, 11).Value = 29729.5
$ws.Cells.Item(31, 12).Value = 3374.3794
$ws.Cells.Item(31, 13).Value = -29434.5
$ws.Cells.Item(31, 14).Value = -3964.3794
$ws.Cells.Item(34, 8).Value = 11955.116
$ws.Cells.Item(34, 9).Value = 29729.5
$ws.Cells.Item(34, 10).Value = 3374.3794
$ws.Cells.Item(34, 11).Value = 29729.5
$ws.Cells.Item(34, 12).Value = 3374.3794
$ws.Cells.Item(34, 13).Value = -29527.5
$ws.Cells.Item(34, 14).Value = -3778.3794
$ws.Cells.Item(132, 8).Value = 17723.5
$ws.Cells.Item(132, 9).Value = 21730.36
$ws.Cells.Item(132, 11).Value = 65191.08
$ws.Cells.Item(132, 13).Value = -62661.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 3946.5
$ws.Cells.Item(63, 10).Value = 5910.8
$ws.Cells.Item(63, 12).Value = 17732.4
$ws.Cells.Item(63, 14).Value = -19230.4
$ws.Cells.Item(64, 8).Value = 1449.0834
$ws.Cells.Item(64, 9).Value = 977.6667
$ws.Cells.Item(64, 10).Value = 2863.3333
$ws.Cells.Item(64, 11).Value = 2933.0001
$ws.Cells.Item(64, 12).Value = 8589.999899999999
$ws.Cells.Item(64, 13).Value = -2663.0001
$ws.Cells.Item(64, 14).Value = -9129.999899999999
$ws.Cells.Item(66, 8).Value = 3946.5
$ws.Cells.Item(66, 10).Value = 5910.8
$ws.Cells.Item(66, 12).Value = 53197.2
$ws.Cells.Item(66, 14).Value = -60685.2
$ws.Cells.Item(67, 8).Value = 1449.0834
$ws.Cells.Item(67, 9).Value = 977.6667
$ws.Cells.Item(67, 10).Value = 2863.3333
$ws.Cells.Item(67, 11).Value = 2933.0001
$ws.Cells.Item(67, 12).Value = 8589.999899999999
$ws.Cells.Item(67, 13).Value = -1997.0001
$ws.Cells.Item(67, 14).Value = -10461.9999
$ws.Cells.Item(68, 8).Value = 4879.222
$ws.Cells.Item(68, 10).Value = 9438.691999999999
$ws.Cells.Item(68, 12).Value = 28316.076
$ws.Cells.Item(68, 14).Value = -29938.076
$ws.Cells.Item(71, 8).Value = 4879.222
$ws.Cells.Item(71, 10).Value = 9438.691999999999
$ws.Cells.Item(71, 12).Value = 84948.22799999999
$ws.Cells.Item(71, 14).Value = -93060.22799999999
$ws.Cells.Item(92, 8).Value = 17862242
$ws.Cells.Item(92, 9).Value = 31250424
$ws.Cells.Item(92, 10).Value = 11333.333
$ws.Cells.Item(92, 11).Value = 93751272
$ws.Cells.Item(92, 12).Value = 33999.999
$ws.Cells.Item(92, 13).Value = -93750024
$ws.Cells.Item(92, 14).Value = -36495.999
$ws.Cells.Item(93, 8).Value = 3000
$ws.Cells.Item(93, 10).Value = 3500
$ws.Cells.Item(93, 12).Value = 10500
$ws.Cells.Item(93, 14).Value = -14244
$ws.Cells.Item(98, 8).Value = 799.5
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 799.5
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 2398.5
$ws.Cells.Item(98, 13).ClearContents()
$ws.Cells.Item(98, 14).Value = -5394.5
$ws.Cells.Item(103, 8).Value = 942.0833
$ws.Cells.Item(103, 9).Value = 571.2857
$ws.Cells.Item(103, 10).Value = 1461.2
$ws.Cells.Item(103, 11).Value = 1713.8571
$ws.Cells.Item(103, 12).Value = 4383.6
$ws.Cells.Item(103, 13).Value = -834.8571000000002
$ws.Cells.Item(103, 14).Value = -6141.6
$ws.Cells.Item(112, 8).Value = 2987.5
$ws.Cells.Item(112, 10).Value = 3666.6667
$ws.Cells.Item(112, 12).Value = 11000.0001
$ws.Cells.Item(112, 14).Value = -13216.0001
$ws.Cells.Item(122, 8).Value = 553
$ws.Cells.Item(122, 10).Value = 623.4
$ws.Cells.Item(122, 12).Value = 5610.599999999999
$ws.Cells.Item(122, 14).Value = -10510.6
$ws.Cells.Item(131, 8).Value = 106087.125
$ws.Cells.Item(131, 10).Value = 110705.68
$ws.Cells.Item(131, 12).Value = 332117.04
$ws.Cells.Item(131, 14).Value = -342197.04

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 2693.6924
$ws.Cells.Item(97, 9).Value = 1356.4445
$ws.Cells.Item(97, 10).Value = 5702.5
$ws.Cells.Item(97, 11).Value = 1356.4445
$ws.Cells.Item(97, 12).Value = 5702.5
$ws.Cells.Item(97, 13).Value = -860.4445000000001
$ws.Cells.Item(97, 14).Value = -6694.5
$ws.Cells.Item(102, 8).Value = 1430.72
$ws.Cells.Item(102, 9).Value = 1374.8096
$ws.Cells.Item(102, 11).Value = 1374.8096
$ws.Cells.Item(102, 13).Value = 247.1904
$ws.Cells.Item(136, 8).Value = 11922.25
$ws.Cells.Item(136, 10).Value = 11922.25
$ws.Cells.Item(136, 12).Value = 35766.75
$ws.Cells.Item(136, 14).Value = -40866.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 3183.2
$ws.Cells.Item(93, 9).Value = 3104
$ws.Cells.Item(93, 11).Value = 3104
$ws.Cells.Item(93, 13).Value = -1856
$ws.Cells.Item(100, 8).Value = 3000
$ws.Cells.Item(100, 9).Value = 3000
$ws.Cells.Item(100, 11).Value = 3000
$ws.Cells.Item(100, 13).Value = -2459

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4700
$ws.Cells.Item(62, 9).Value = 4000
$ws.Cells.Item(62, 10).Value = 5400
$ws.Cells.Item(62, 11).Value = 4000
$ws.Cells.Item(62, 12).Value = 5400
$ws.Cells.Item(62, 13).Value = -3376
$ws.Cells.Item(62, 14).Value = -6648
$ws.Cells.Item(65, 8).Value = 4700
$ws.Cells.Item(65, 9).Value = 4000
$ws.Cells.Item(65, 10).Value = 5400
$ws.Cells.Item(65, 11).Value = 20000
$ws.Cells.Item(65, 12).Value = 27000
$ws.Cells.Item(65, 13).Value = -16880
$ws.Cells.Item(65, 14).Value = -33240
$ws.Cells.Item(102, 8).Value = 36333.332
$ws.Cells.Item(102, 10).Value = 36333.332
$ws.Cells.Item(102, 12).Value = 36333.332
$ws.Cells.Item(102, 14).Value = -42823.332
$ws.Cells.Item(132, 8).Value = 2192.926
$ws.Cells.Item(132, 9).Value = 2100.6875
$ws.Cells.Item(132, 10).Value = 2327.0908
$ws.Cells.Item(132, 11).Value = 6302.0625
$ws.Cells.Item(132, 12).Value = 6981.2724
$ws.Cells.Item(132, 13).Value = -3772.0625
$ws.Cells.Item(132, 14).Value = -12041.2724
